$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$values = @(
    "51-23=",
    "12+66=",
    "21+56=",
    "42-16=",
    "85-23=",
    "15+71=",
    "53+0=",
    "20+2=",
    "96-1=",
    "26+16=",
    "28+16=",
    "92-92=",
    "43-39=",
    "77+15=",
    "39-38=",
    "10+41=",
    "10+52=",
    "91-13=",
    "64+17=",
    "59+21=",
    "31+25=",
    "48-15=",
    "8+17=",
    "56-11=",
    "41+7=",
    "68+20=",
    "37+61=",
    "77+14=",
    "83-51=",
    "72-12=",
    "2+63=",
    "99-0=",
    "7+62=",
    "83+6=",
    "45+46=",
    "0+97=",
    "21-2=",
    "68-20=",
    "52-50=",
    "65-23=",
    "82-19=",
    "16+42=",
    "61-21=",
    "18+58=",
    "3+16=",
    "38+29=",
    "28-13=",
    "68-17=",
    "77-22=",
    "24+30=",
    "67-1=",
    "25+57=",
    "69+2=",
    "49-24=",
    "59+4=",
    "76+0=",
    "31+58=",
    "55-13=",
    "49+32=",
    "73-55=",
    "13+39=",
    "96-51=",
    "39-19=",
    "70-30=",
    "81-61=",
    "0+60=",
    "60-25=",
    "58-0=",
    "87-85=",
    "77-9=",
    "53-23=",
    "97-15=",
    "43+33=",
    "84+8=",
    "30+42=",
    "26+31=",
    "84-45=",
    "42-38=",
    "92-75=",
    "54-20=",
    "57-6=",
    "25+11=",
    "98-27=",
    "70-15=",
    "15-6=",
    "46-45=",
    "30+20=",
    "62-6=",
    "48-13=",
    "26+70=",
    "56-14=",
    "43-11=",
    "16+41=",
    "72-51=",
    "65-3=",
    "3+53=",
    "50-37=",
    "15+62=",
    "45-41=",
    "2+68="
)

$idx = 0
for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $t.Cell($r, $c).Range.Text = $values[$idx]
        $idx++
    }
}
Write-Host "Updated" $idx "cells"
